$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently ends at row 8:
#   A8/B8/C8/D8 = "13"/"14"/"14"/"6" (shared strings), E8 = 15.0, F8 = 16.0
#
# Target state:
#   Row 8 gets its E/F values swapped -> E8 = 16.0, F8 = 15.0
#   A brand-new row 9 is appended, identical in shape/style/strings to row 8,
#   also with E9 = 16.0, F9 = 15.0

# 1) Duplicate row 8 (values, shared-string refs, and styles) down into row 9
#    using Range.Copy(Destination) so formatting/shared strings survive intact.
$ws.Range("A8:F8").Copy($ws.Range("A9"))

# 2) Swap E8/F8 on the original row.
$ws.Range("E8").Value = 16.0
$ws.Range("F8").Value = 15.0

# 3) Swap E9/F9 on the newly duplicated row (it was copied with the old
#    15.0/16.0 values from row 8, so apply the same swap here).
$ws.Range("E9").Value = 16.0
$ws.Range("F9").Value = 15.0
